$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run SGNN dialog act annotations: update DAMSLTag (I) and DialogAct (J) columns
# for the rows identified by the re-annotation pass.
$ws.Range("I3").Value = 'sd'
$ws.Range("J3").Value = 'Statement-non-opinion'
$ws.Range("I25").Value = '%'
$ws.Range("J25").Value = 'Uninterpretable'
$ws.Range("I26").Value = 'sd'
$ws.Range("J26").Value = 'Statement-non-opinion'
$ws.Range("I40").Value = 'sd'
$ws.Range("J40").Value = 'Statement-non-opinion'
$ws.Range("I41").Value = 'b'
$ws.Range("J41").Value = 'Acknowledge (Backchannel)'
$ws.Range("I45").Value = 'b'
$ws.Range("J45").Value = 'Acknowledge (Backchannel)'
$ws.Range("I67").Value = 'ba'
$ws.Range("J67").Value = 'Appreciation'
$ws.Range("I69").Value = 'b'
$ws.Range("J69").Value = 'Acknowledge (Backchannel)'
$ws.Range("I91").Value = 'sd'
$ws.Range("J91").Value = 'Statement-non-opinion'
$ws.Range("I94").Value = 'b'
$ws.Range("J94").Value = 'Acknowledge (Backchannel)'
$ws.Range("I108").Value = 'sd'
$ws.Range("J108").Value = 'Statement-non-opinion'
$ws.Range("I129").Value = 'sv'
$ws.Range("J129").Value = 'Statement-opinion'
$ws.Range("I136").Value = 'aa'
$ws.Range("J136").Value = 'Agree/Accept'
$ws.Range("I148").Value = 'sd'
$ws.Range("J148").Value = 'Statement-non-opinion'
$ws.Range("I161").Value = 'aa'
$ws.Range("J161").Value = 'Agree/Accept'
$ws.Range("I171").Value = '%'
$ws.Range("J171").Value = 'Uninterpretable'
$ws.Range("I172").Value = '%'
$ws.Range("J172").Value = 'Uninterpretable'
$ws.Range("I173").Value = 'sd'
$ws.Range("J173").Value = 'Statement-non-opinion'
$ws.Range("I184").Value = 'aa'
$ws.Range("J184").Value = 'Agree/Accept'
$ws.Range("I188").Value = 'aa'
$ws.Range("J188").Value = 'Agree/Accept'
$ws.Range("I189").Value = '%'
$ws.Range("J189").Value = 'Uninterpretable'
$ws.Range("I190").Value = 'sv'
$ws.Range("J190").Value = 'Statement-opinion'
$ws.Range("I198").Value = 'b'
$ws.Range("J198").Value = 'Acknowledge (Backchannel)'
$ws.Range("I202").Value = '%'
$ws.Range("J202").Value = 'Uninterpretable'
$ws.Range("I204").Value = 'sd'
$ws.Range("J204").Value = 'Statement-non-opinion'
$ws.Range("I207").Value = 'aa'
$ws.Range("J207").Value = 'Agree/Accept'
$ws.Range("I209").Value = 'b'
$ws.Range("J209").Value = 'Acknowledge (Backchannel)'
$ws.Range("I211").Value = 'sv'
$ws.Range("J211").Value = 'Statement-opinion'
$ws.Range("I225").Value = 'sd'
$ws.Range("J225").Value = 'Statement-non-opinion'
$ws.Range("I226").Value = 'sv'
$ws.Range("J226").Value = 'Statement-opinion'
$ws.Range("I232").Value = 'ba'
$ws.Range("J232").Value = 'Appreciation'
$ws.Range("I242").Value = 'sv'
$ws.Range("J242").Value = 'Statement-opinion'
$ws.Range("I253").Value = 'ba'
$ws.Range("J253").Value = 'Appreciation'
$ws.Range("I260").Value = 'aa'
$ws.Range("J260").Value = 'Agree/Accept'
$ws.Range("I268").Value = 'aa'
$ws.Range("J268").Value = 'Agree/Accept'
$ws.Range("I270").Value = 'sd'
$ws.Range("J270").Value = 'Statement-non-opinion'
$ws.Range("I272").Value = 'ba'
$ws.Range("J272").Value = 'Appreciation'
$ws.Range("I274").Value = 'ba'
$ws.Range("J274").Value = 'Appreciation'
$ws.Range("I277").Value = 'sv'
$ws.Range("J277").Value = 'Statement-opinion'
$ws.Range("I280").Value = 'ba'
$ws.Range("J280").Value = 'Appreciation'
$ws.Range("I281").Value = 'b'
$ws.Range("J281").Value = 'Acknowledge (Backchannel)'
$ws.Range("I285").Value = 'b'
$ws.Range("J285").Value = 'Acknowledge (Backchannel)'
$ws.Range("I295").Value = 'sv'
$ws.Range("J295").Value = 'Statement-opinion'
$ws.Range("I296").Value = 'sv'
$ws.Range("J296").Value = 'Statement-opinion'
$ws.Range("I306").Value = '%'
$ws.Range("J306").Value = 'Uninterpretable'
$ws.Range("I309").Value = 'sd'
$ws.Range("J309").Value = 'Statement-non-opinion'
$ws.Range("I312").Value = 'aa'
$ws.Range("J312").Value = 'Agree/Accept'
$ws.Range("I313").Value = 'sd'
$ws.Range("J313").Value = 'Statement-non-opinion'
$ws.Range("I316").Value = 'ba'
$ws.Range("J316").Value = 'Appreciation'
$ws.Range("I319").Value = 'ba'
$ws.Range("J319").Value = 'Appreciation'
$ws.Range("I340").Value = 'sv'
$ws.Range("J340").Value = 'Statement-opinion'
$ws.Range("I345").Value = 'sd'
$ws.Range("J345").Value = 'Statement-non-opinion'
$ws.Range("I370").Value = 'b'
$ws.Range("J370").Value = 'Acknowledge (Backchannel)'
$ws.Range("I380").Value = 'sd'
$ws.Range("J380").Value = 'Statement-non-opinion'
$ws.Range("I395").Value = 'sd'
$ws.Range("J395").Value = 'Statement-non-opinion'
$ws.Range("I397").Value = 'b'
$ws.Range("J397").Value = 'Acknowledge (Backchannel)'
$ws.Range("I435").Value = '%'
$ws.Range("J435").Value = 'Uninterpretable'
$ws.Range("I441").Value = 'qy'
$ws.Range("J441").Value = 'Yes-No-Question'
$ws.Range("I444").Value = 'sv'
$ws.Range("J444").Value = 'Statement-opinion'
$ws.Range("I445").Value = 'sd'
$ws.Range("J445").Value = 'Statement-non-opinion'
$ws.Range("I457").Value = 'sv'
$ws.Range("J457").Value = 'Statement-opinion'
$ws.Range("I484").Value = 'ba'
$ws.Range("J484").Value = 'Appreciation'
$ws.Range("I485").Value = 'b'
$ws.Range("J485").Value = 'Acknowledge (Backchannel)'
$ws.Range("I487").Value = 'ba'
$ws.Range("J487").Value = 'Appreciation'
$ws.Range("I489").Value = 'sv'
$ws.Range("J489").Value = 'Statement-opinion'
$ws.Range("I490").Value = 'sd'
$ws.Range("J490").Value = 'Statement-non-opinion'
$ws.Range("I495").Value = 'sv'
$ws.Range("J495").Value = 'Statement-opinion'
$ws.Range("I501").Value = 'sv'
$ws.Range("J501").Value = 'Statement-opinion'
$ws.Range("I505").Value = 'sd'
$ws.Range("J505").Value = 'Statement-non-opinion'
$ws.Range("I516").Value = 'sd'
$ws.Range("J516").Value = 'Statement-non-opinion'
$ws.Range("I526").Value = 'sv'
$ws.Range("J526").Value = 'Statement-opinion'
$ws.Range("I541").Value = 'sd'
$ws.Range("J541").Value = 'Statement-non-opinion'
$ws.Range("I543").Value = 'b'
$ws.Range("J543").Value = 'Acknowledge (Backchannel)'
$ws.Range("I552").Value = 'sv'
$ws.Range("J552").Value = 'Statement-opinion'
$ws.Range("I557").Value = 'sv'
$ws.Range("J557").Value = 'Statement-opinion'
$ws.Range("I558").Value = 'sv'
$ws.Range("J558").Value = 'Statement-opinion'
$ws.Range("I577").Value = 'b'
$ws.Range("J577").Value = 'Acknowledge (Backchannel)'
$ws.Range("I578").Value = 'sd'
$ws.Range("J578").Value = 'Statement-non-opinion'
$ws.Range("I580").Value = 'sv'
$ws.Range("J580").Value = 'Statement-opinion'
$ws.Range("I604").Value = 'sd'
$ws.Range("J604").Value = 'Statement-non-opinion'
$ws.Range("I612").Value = 'aa'
$ws.Range("J612").Value = 'Agree/Accept'
$ws.Range("I613").Value = 'sv'
$ws.Range("J613").Value = 'Statement-opinion'
$ws.Range("I627").Value = 'sv'
$ws.Range("J627").Value = 'Statement-opinion'
